# Commit from Develop branch
#
# Insert a new worksheet "DifferentActions" right after the first sheet
# (CreateDeleteCustFunctionality) and before "OpenAccountFunctionality",
# populate it with the new test-case rows, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)

# Add the new worksheet immediately after the first sheet.
$newSheet = $wb.Worksheets.Add($null, $firstSheet)
$newSheet.Name = "DifferentActions"

# Populate the new sheet. The write order below mirrors how the data was
# originally authored (TC Name / Status / RunMode columns first, then the
# Window_Authorization row, then the UserName/Password columns, etc.) so the
# shared-string table is rebuilt in the same order as the target workbook.
$newSheet.Range("A1").Value = "TC Name"
$newSheet.Range("D1").Value = "Status"
$newSheet.Range("E1").Value = "RunMode"

$newSheet.Range("A2").Value = "Window_Authorization"

$newSheet.Range("B1").Value = "UserName"
$newSheet.Range("C1").Value = "Password"

$newSheet.Range("B2").Value = "admin"
$newSheet.Range("C2").Value = "admin"
$newSheet.Range("E2").Value = "Y"

$newSheet.Range("A3").Value = "Upload_File_Using_Send_Keys"
$newSheet.Range("E3").Value = "Y"

# Restore the selection on the original first sheet before switching away.
$firstSheet.Activate()
$null = $firstSheet.Range("G1").Select()

# Make the newly inserted sheet the active / selected tab, matching the
# selection left on it after data entry.
$newSheet.Activate()
$null = $newSheet.Range("E3").Select()
